$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4964252412319183
$ws.Range("B1").Value = 2.499825000762939
$ws.Range("C1").Value = 6.400107860565186
$ws.Range("D1").Value = 1.557233095169067
$ws.Range("E1").Value = 0.8968315720558167
